$wb = $excel.ActiveWorkbook

# Add a new "Comments" column (column E) header to the four history sheets.
$wsWithdraw = $wb.Worksheets("Withdraw History")
$wsDeposit  = $wb.Worksheets("Deposit History")
$wsTransfer = $wb.Worksheets("Transfer History")
$wsAbsolute = $wb.Worksheets("Absolute History")

$wsWithdraw.Range("E1").Value = "Comments"
$wsDeposit.Range("E1").Value  = "Comments"
$wsTransfer.Range("E1").Value = "Comments"
$wsAbsolute.Range("E1").Value = "Comments"

# Restore each sheet's own selection/active-cell state.
$wsDeposit.Range("E1").Select()
$wsTransfer.Range("E1").Select()
$wsAbsolute.Range("E5").Select()

# Make "Withdraw History" the active sheet/tab, with E1 selected, last
# so it becomes the workbook's active tab.
$wsWithdraw.Activate()
$wsWithdraw.Range("E1").Select()
